# Refresh the team-specific transition-probability matrix on Sheet1 with the
# updated simulation results ("added more games, sped up simulate game
# logic, and drafted optimization logic"). Each row is a probability
# distribution over outcome columns (B:S); only the cells whose underlying
# counts changed in the new simulation run are touched here.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1201550387596899
$ws.Range("C2").Value = 0.6162790697674418
$ws.Range("J2").Value = 0.01937984496124031
$ws.Range("P2").Value = 0.1511627906976744
$ws.Range("S2").Value = 0.09302325581395349

# Row 3
$ws.Range("B3").Value = 0.01129943502824859
$ws.Range("C3").Value = 0.02259887005649718
$ws.Range("J3").Value = 0.05084745762711865
$ws.Range("P3").Value = 0.7401129943502824
$ws.Range("S3").Value = 0.1751412429378531

# Row 4
$ws.Range("J4").Value = 0.1578947368421053
$ws.Range("P4").Value = 0.6842105263157895
$ws.Range("S4").Value = 0.1578947368421053

# Row 6
$ws.Range("B6").Value = 0.06217616580310881
$ws.Range("D6").Value = 0.02072538860103627
$ws.Range("F6").Value = 0.05181347150259067
$ws.Range("J6").Value = 0.2694300518134715
$ws.Range("O6").Value = 0.0155440414507772
$ws.Range("Q6").Value = 0.1658031088082902
$ws.Range("R6").Value = 0.04663212435233161
$ws.Range("S6").Value = 0.3678756476683938

# Row 7
$ws.Range("B7").Value = 0.1136363636363636
$ws.Range("F7").Value = 0.04545454545454546
$ws.Range("J7").Value = 0.0625
$ws.Range("O7").Value = 0.1136363636363636
$ws.Range("Q7").Value = 0.1477272727272727
$ws.Range("R7").Value = 0.1022727272727273
$ws.Range("S7").Value = 0.4147727272727273

# Row 8
$ws.Range("B8").Value = 0.1316397228637413
$ws.Range("D8").Value = 0.01847575057736721
$ws.Range("F8").Value = 0.06928406466512702
$ws.Range("J8").Value = 0.1270207852193996
$ws.Range("O8").Value = 0.01616628175519631
$ws.Range("Q8").Value = 0.1801385681293303
$ws.Range("R8").Value = 0.06004618937644342
$ws.Range("S8").Value = 0.3972286374133949

# Row 9
$ws.Range("B9").Value = 0.1049382716049383
$ws.Range("D9").Value = 0.02469135802469136
$ws.Range("F9").Value = 0.06172839506172839
$ws.Range("J9").Value = 0.08641975308641975
$ws.Range("O9").Value = 0.006172839506172839
$ws.Range("Q9").Value = 0.1790123456790123
$ws.Range("R9").Value = 0.06172839506172839
$ws.Range("S9").Value = 0.4753086419753086

# Row 10
$ws.Range("B10").Value = 0.08952702702702703
$ws.Range("D10").Value = 0.01942567567567568
$ws.Range("E10").Value = 0.0008445945945945946
$ws.Range("F10").Value = 0.0633445945945946
$ws.Range("J10").Value = 0.183277027027027
$ws.Range("O10").Value = 0.02449324324324324
$ws.Range("Q10").Value = 0.1875
$ws.Range("R10").Value = 0.08699324324324324
$ws.Range("S10").Value = 0.3445945945945946

# Row 11
$ws.Range("G11").Value = 0.1440677966101695
$ws.Range("J11").Value = 0.1059322033898305
$ws.Range("K11").Value = 0.1313559322033898
$ws.Range("L11").Value = 0.6186440677966102

# Row 12
$ws.Range("G12").Value = 0.7278911564625851
$ws.Range("J12").Value = 0.217687074829932
$ws.Range("L12").Value = 0.0272108843537415
$ws.Range("S12").Value = 0.0272108843537415

# Row 13
$ws.Range("F13").Value = 0.02040816326530612
$ws.Range("G13").Value = 0.7142857142857143
$ws.Range("J13").Value = 0.2448979591836735
$ws.Range("S13").Value = 0.02040816326530612

# Row 15
$ws.Range("F15").Value = 0.02173913043478261
$ws.Range("H15").Value = 0.1739130434782609
$ws.Range("I15").Value = 0.04347826086956522
$ws.Range("J15").Value = 0.358695652173913
$ws.Range("K15").Value = 0.07065217391304347
$ws.Range("M15").Value = 0.01630434782608696
$ws.Range("O15").Value = 0.04347826086956522
$ws.Range("S15").Value = 0.2717391304347826

# Row 16
$ws.Range("F16").Value = 0.08947368421052632
$ws.Range("H16").Value = 0.2052631578947368
$ws.Range("I16").Value = 0.08947368421052632
$ws.Range("J16").Value = 0.3736842105263158
$ws.Range("K16").Value = 0.07368421052631578
$ws.Range("M16").Value = 0.02631578947368421
$ws.Range("O16").Value = 0.02105263157894737
$ws.Range("S16").Value = 0.1210526315789474

# Row 17
$ws.Range("F17").Value = 0.01025641025641026
$ws.Range("H17").Value = 0.182051282051282
$ws.Range("I17").Value = 0.1025641025641026
$ws.Range("J17").Value = 0.3974358974358974
$ws.Range("K17").Value = 0.09487179487179487
$ws.Range("M17").Value = 0.02820512820512821
$ws.Range("O17").Value = 0.05641025641025641
$ws.Range("S17").Value = 0.1282051282051282

# Row 18
$ws.Range("F18").Value = 0.02395209580838323
$ws.Range("H18").Value = 0.1377245508982036
$ws.Range("I18").Value = 0.07784431137724551
$ws.Range("J18").Value = 0.4311377245508982
$ws.Range("K18").Value = 0.1437125748502994
$ws.Range("M18").Value = 0.03592814371257485
$ws.Range("O18").Value = 0.04790419161676647
$ws.Range("S18").Value = 0.1017964071856287

# Row 19
$ws.Range("F19").Value = 0.009777777777777778
$ws.Range("H19").Value = 0.232
$ws.Range("I19").Value = 0.07644444444444444
$ws.Range("J19").Value = 0.3733333333333334
$ws.Range("K19").Value = 0.1075555555555556
$ws.Range("M19").Value = 0.02222222222222222
$ws.Range("N19").Value = 0.0008888888888888889
$ws.Range("O19").Value = 0.06044444444444445
$ws.Range("S19").Value = 0.1173333333333333
